$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must stay as literal text
# (matching the source inlineStr cells), so force text format, assign, then
# restore the default "Normal" style so no stray formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.013.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.025.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.641"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.08"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +14.25%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.85"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.372"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0746"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.85%  "

$ws.Range("E12").Value = "  -1.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.894"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.326.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("E16").Value = "  +4.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +18.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.031.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.000.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0866"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +24.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.120"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.83%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.31%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.45%  "

$ws.Range("E33").Value = "  +24.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.95%  "

$ws.Range("E35").Value = "  +3.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.64%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +21.92%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.104"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +16.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.20%  "

$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.75%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +22.53%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0216"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.02%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.419.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.96%  "
